$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5781
$ws.Range("L3").Value = 6307
$ws.Range("E4").Value = 1804
$ws.Range("F4").Value = 1714
$ws.Range("L4").Value = 1562
$ws.Range("L5").Value = 377
$ws.Range("L6").Value = 5179
$ws.Range("E7").Value = 22830
$ws.Range("F7").Value = 21038
$ws.Range("L7").Value = 19206

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 212

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L4").Value = 89
$ws.Range("L5").Value = 42
$ws.Range("L6").Value = 311
$ws.Range("L7").Value = 1265

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 219
$ws.Range("L3").Value = 256
$ws.Range("L7").Value = 732

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 369

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 136
$ws.Range("L7").Value = 335

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L5").Value = 69
$ws.Range("L7").Value = 614
$ws.Range("L8").Value = 1265
$ws.Range("L11").Value = 317
$ws.Range("L19").Value = 522
$ws.Range("L20").Value = 479
$ws.Range("L23").Value = 210
$ws.Range("L24").Value = 54
$ws.Range("L29").Value = 1082
$ws.Range("L31").Value = 188
$ws.Range("L37").Value = 732
$ws.Range("L39").Value = 12
$ws.Range("L42").Value = 620
$ws.Range("L44").Value = 134
$ws.Range("L47").Value = 132
$ws.Range("L49").Value = 104
$ws.Range("L50").Value = 94
$ws.Range("L51").Value = 240
$ws.Range("L52").Value = 400
$ws.Range("L53").Value = 212
$ws.Range("L55").Value = 198
$ws.Range("L60").Value = 123
$ws.Range("E63").Value = 323
$ws.Range("F63").Value = 191
$ws.Range("L63").Value = 61
$ws.Range("L65").Value = 369
$ws.Range("L78").Value = 244
$ws.Range("L79").Value = 530
$ws.Range("L82").Value = 27
$ws.Range("L85").Value = 951
$ws.Range("L90").Value = 198
$ws.Range("L94").Value = 236
$ws.Range("L97").Value = 157
$ws.Range("L99").Value = 335
$ws.Range("E101").Value = 22830
$ws.Range("F101").Value = 21038
$ws.Range("L101").Value = 19206

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 417
$ws.Range("L7").Value = 1082

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L6").Value = 142
$ws.Range("L7").Value = 522

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L6").Value = 167
$ws.Range("L7").Value = 620

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 244

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 198

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 85
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 170
$ws.Range("L3").Value = 169
$ws.Range("L7").Value = 530

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 166
$ws.Range("L6").Value = 116
$ws.Range("L7").Value = 479

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 196
$ws.Range("L7").Value = 614

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 55
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 51
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L2").Value = 5
$ws.Range("L6").Value = 12

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 317

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 157

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 198

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L4").Value = 35
$ws.Range("L7").Value = 240

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 392
$ws.Range("L6").Value = 197
$ws.Range("L7").Value = 951

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L4").Value = 28
$ws.Range("L5").Value = 12
$ws.Range("L7").Value = 400
